$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "DDL" column (E), shifting
# it (and its data) one column to the right, then populate the new
# column with the "ObjectType" header/value.
$ws.Range("E1").EntireColumn.Insert()

$ws.Range("E1").Value = "ObjectType"
$ws.Range("E2").Value = "Sconosciuto"
